# Reservation Page Controller: when a customer confirms their order, the
# confirmation method now saves the customer information (first name, last
# name, email/username, phone, and confirmation ID) into the Customers
# sheet for each newly booked room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# Row 21 -> Room 20's reservation confirmation
$ws.Range("B21").Value = "Nathan"
$ws.Range("C21").Value = "Wahba"
$ws.Range("D21").Value = "nwahba02"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "123"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = "9iOO1"

# Row 22 -> Room 21's reservation confirmation
$ws.Range("B22").Value = "Nathan"
$ws.Range("C22").Value = "Wahba"
$ws.Range("D22").Value = "nwahba02"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "123"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Value = "kB0c7"
